# Updated symbol list (Price / Volume(1h) columns) to refresh cryptocurrency
# quote data, matching the automated "Updated symbol list" GitHub Actions commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text ("@") number format on the cells we touch so the refreshed
# Price/Volume values are written back as plain text (e.g. "2.51%"),
# matching how the source data is stored, instead of being auto-parsed
# into numeric/percentage values by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

$ws.Range("D2").Value = "294.88"
$ws.Range("E2").Value = "2.51%"
$ws.Range("D3").Value = "31.09"
$ws.Range("E3").Value = "0.81%"
$ws.Range("D4").Value = "4.933"
$ws.Range("E4").Value = "0.81%"
$ws.Range("D5").Value = "0.07385"
$ws.Range("E5").Value = "4.02%"
$ws.Range("D6").Value = "2.306"
$ws.Range("E6").Value = "30.08%"
$ws.Range("D7").Value = "7.687"
$ws.Range("E7").Value = "0.67%"
$ws.Range("D8").Value = "3.753"
$ws.Range("D9").Value = "0.9141"
$ws.Range("E9").Value = "2.00%"
$ws.Range("D10").Value = "0.1692"
$ws.Range("D11").Value = "0.08338"
$ws.Range("E11").Value = "11.04%"
$ws.Range("D12").Value = "0.08300"
$ws.Range("E12").Value = "3.80%"
$ws.Range("D13").Value = "0.03116"
$ws.Range("E13").Value = "3.90%"
$ws.Range("E14").Value = "0.86%"
$ws.Range("D15").Value = "0.001511"
$ws.Range("E15").Value = "0.85%"
$ws.Range("D16").Value = "0.005698"
$ws.Range("E16").Value = "1.03%"
$ws.Range("E17").Value = "0.40%"
$ws.Range("E18").Value = "-1.36%"
$ws.Range("E19").Value = "1.65%"
$ws.Range("D20").Value = "0.1303"
$ws.Range("D21").Value = "3.968"
$ws.Range("E21").Value = "-2.23%"
$ws.Range("D22").Value = "0.2100"
$ws.Range("E22").Value = "4.70%"
$ws.Range("D23").Value = "0.04544"
$ws.Range("E23").Value = "1.15%"
$ws.Range("D24").Value = "0.001209"
$ws.Range("E24").Value = "-0.44%"
$ws.Range("D25").Value = "0.004334"
$ws.Range("E25").Value = "-6.51%"
$ws.Range("D26").Value = "0.0001301"
$ws.Range("E26").Value = "3.78%"
$ws.Range("D39").Value = "0.01606"
$ws.Range("E39").Value = "-0.31%"
$ws.Range("D40").Value = "0.04457"
$ws.Range("E40").Value = "2.78%"
$ws.Range("D41").Value = "0.007337"
$ws.Range("E41").Value = "-0.89%"
$ws.Range("D43").Value = "0.1326"
$ws.Range("E43").Value = "1.86%"
$ws.Range("D44").Value = "0.002061"
$ws.Range("E44").Value = "2.72%"
$ws.Range("D45").Value = "0.009189"
$ws.Range("E45").Value = "-11.13%"
$ws.Range("D46").Value = "0.00006022"
$ws.Range("E46").Value = "2.15%"
$ws.Range("E47").Value = "-0.20%"
$ws.Range("E48").Value = "0.93%"
$ws.Range("E49").Value = "-3.62%"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").Value = "-0.20%"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").Value = "-0.20%"
